$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A (the plain numeric index column 0..5) and shift the
# States / Capitals / Population columns one to the left.
$ws.Columns("A").Delete()
